# Auto-generated edit script: applies cell-value updates to match target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1081.0769
$ws.Range("J2").Value = 1444.7142
$ws.Range("L2").Value = 1444.7142
$ws.Range("N2").Value = -1670.7142
$ws.Range("H17").Value = 864.4754
$ws.Range("J17").Value = 870.0678
$ws.Range("L17").Value = 2610.2034
$ws.Range("N17").Value = -2946.2034
$ws.Range("H32").Value = 11691.0625
$ws.Range("J32").Value = 10106.6
$ws.Range("L32").Value = 10106.6
$ws.Range("N32").Value = -10758.6
$ws.Range("H33").Value = 438.21875
$ws.Range("I33").Value = 144.36363
$ws.Range("K33").Value = 144.36363
$ws.Range("M33").Value = 84.63637
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H116").Value = 15709.066
$ws.Range("I116").Value = 17745.072
$ws.Range("J116").Value = 13927.5625
$ws.Range("K116").Value = 17745.072
$ws.Range("L116").Value = 13927.5625
$ws.Range("M116").Value = -14303.072
$ws.Range("N116").Value = -20811.5625
$ws.Range("H138").Value = 41927.96
$ws.Range("I138").Value = 1912.3043
$ws.Range("K138").Value = 5736.9129
$ws.Range("M138").Value = -596.9129000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19128.121
$ws.Range("I32").Value = 21122.46
$ws.Range("J32").Value = 1843.8334
$ws.Range("K32").Value = 21122.46
$ws.Range("L32").Value = 1843.8334
$ws.Range("M32").Value = -20835.46
$ws.Range("N32").Value = -2417.8334
$ws.Range("H122").Value = 1969.2142
$ws.Range("I122").Value = 2047.5
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 6142.5
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -3692.5
$ws.Range("N122").Value = -9398.5
$ws.Range("H132").Value = 1653.2642
$ws.Range("I132").Value = 1075.3103
$ws.Range("K132").Value = 3225.9309
$ws.Range("M132").Value = -695.9309000000003
$ws.Range("H139").Value = 141663.33
$ws.Range("J139").Value = 141663.33
$ws.Range("L139").Value = 141663.33
$ws.Range("N139").Value = -151943.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 291.5
$ws.Range("I29").Value = 291.5
$ws.Range("K29").Value = 291.5
$ws.Range("M29").Value = -2.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2688.5
$ws.Range("J16").Value = 3185
$ws.Range("L16").Value = 3185
$ws.Range("N16").Value = -3759
$ws.Range("H31").Value = 3227649.8
$ws.Range("I31").Value = 3847159
$ws.Range("J31").Value = 6202.4
$ws.Range("K31").Value = 3847159
$ws.Range("L31").Value = 6202.4
$ws.Range("M31").Value = -3846864
$ws.Range("N31").Value = -6792.4
$ws.Range("H34").Value = 3227649.8
$ws.Range("I34").Value = 3847159
$ws.Range("J34").Value = 6202.4
$ws.Range("K34").Value = 3847159
$ws.Range("L34").Value = 6202.4
$ws.Range("M34").Value = -3846957
$ws.Range("N34").Value = -6606.4
$ws.Range("H58").Value = 13735.692
$ws.Range("I58").Value = 1409.5834
$ws.Range("K58").Value = 1409.5834
$ws.Range("M58").Value = -1206.5834
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -112488
$ws.Range("H81").Value = 35000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 35000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H99").Value = 7064.857
$ws.Range("J99").Value = 9784.799999999999
$ws.Range("L99").Value = 9784.799999999999
$ws.Range("N99").Value = -12780.8
$ws.Range("H113").Value = 2688.5
$ws.Range("J113").Value = 3185
$ws.Range("L113").Value = 3185
$ws.Range("N113").Value = -7525
$ws.Range("H126").Value = 7064.857
$ws.Range("J126").Value = 9784.799999999999
$ws.Range("L126").Value = 29354.4
$ws.Range("N126").Value = -34294.39999999999
$ws.Range("H136").Value = 13735.692
$ws.Range("I136").Value = 1409.5834
$ws.Range("K136").Value = 4228.7502
$ws.Range("M136").Value = -1678.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I38").Value = 17.8
$ws.Range("J38").Value = 166666690
$ws.Range("K38").Value = 53.40000000000001
$ws.Range("L38").Value = 500000070
$ws.Range("M38").Value = 293.6
$ws.Range("N38").Value = -500000764
$ws.Range("H64").Value = 4833.125
$ws.Range("J64").Value = 4999.8184
$ws.Range("L64").Value = 14999.4552
$ws.Range("N64").Value = -15539.4552
$ws.Range("H67").Value = 4833.125
$ws.Range("J67").Value = 4999.8184
$ws.Range("L67").Value = 14999.4552
$ws.Range("N67").Value = -16871.4552
$ws.Range("H140").Value = 4383.5557
$ws.Range("I140").Value = 4383.5557
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 13150.6671
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -7970.667099999999
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2873.2307
$ws.Range("J7").Value = 3049.6667
$ws.Range("L7").Value = 3049.6667
$ws.Range("N7").Value = -3273.6667
$ws.Range("H40").Value = 1913.4231
$ws.Range("I40").Value = 1789.96
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 1789.96
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -1653.96
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 5088.9165
$ws.Range("H126").Value = 2873.2307
$ws.Range("J126").Value = 3049.6667
$ws.Range("L126").Value = 9149.000100000001
$ws.Range("N126").Value = -14089.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7798.35
$ws.Range("J81").Value = 5352.7144
$ws.Range("L81").Value = 10705.4288
$ws.Range("N81").Value = -12827.4288
$ws.Range("H84").Value = 7798.35
$ws.Range("J84").Value = 5352.7144
$ws.Range("L84").Value = 53527.144
$ws.Range("N84").Value = -64135.144
